$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Legacy VML picture shapes (title-page logos): Word renumbers the
#    o:spid values and marks the wrap as "edited" when it resaves a
#    document containing these drawings. Reproduce that by replacing
#    the paragraph that hosts the two <w:pict> drawings with an
#    equivalent paragraph whose v:shape attributes match.
# ------------------------------------------------------------------
$picParaIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $pxml = $p.Range.WordOpenXML
    if ($pxml -like "*5300_IBMpos*") {
        $picParaIdx = $i
        break
    }
}
if ($picParaIdx -eq -1) { throw "Could not locate the picture paragraph" }

$picPara = $d.Paragraphs.Item($picParaIdx)
$picRng = $d.Range($picPara.Range.Start, $picPara.Range.End)

$picXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="46E1C894" w14:textId="77777777" w:rsidR="00322587" w:rsidRDefault="00000000"><w:pPr><w:framePr w:w="5580" w:h="12430" w:hRule="exact" w:hSpace="187" w:wrap="around" w:vAnchor="page" w:hAnchor="page" w:y="2026" w:anchorLock="1"/><w:shd w:val="solid" w:color="FFFFFF" w:fill="FFFFFF"/></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:pict w14:anchorId="7F723D40"><v:shapetype id="_x0000_t75" coordsize="21600,21600" o:spt="75" o:preferrelative="t" path="m@4@5l@4@11@9@11@9@5xe" filled="f" stroked="f"><v:stroke joinstyle="miter"/><v:formulas><v:f eqn="if lineDrawn pixelLineWidth 0"/><v:f eqn="sum @0 1 0"/><v:f eqn="sum 0 0 @1"/><v:f eqn="prod @2 1 2"/><v:f eqn="prod @3 21600 pixelWidth"/><v:f eqn="prod @3 21600 pixelHeight"/><v:f eqn="sum @0 0 1"/><v:f eqn="prod @6 1 2"/><v:f eqn="prod @7 21600 pixelWidth"/><v:f eqn="sum @8 21600 0"/><v:f eqn="prod @7 21600 pixelHeight"/><v:f eqn="sum @10 21600 0"/></v:formulas><v:path o:extrusionok="f" gradientshapeok="t" o:connecttype="rect"/><o:lock v:ext="edit" aspectratio="t"/></v:shapetype><v:shape id="Picture 1" o:spid="_x0000_s1027" type="#_x0000_t75" alt="Description: Description: 5300_IBMpos" style="position:absolute;margin-left:98.95pt;margin-top:193.35pt;width:64.5pt;height:23.25pt;z-index:2;visibility:visible;mso-wrap-style:square;mso-wrap-edited:f;mso-width-percent:0;mso-height-percent:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:page;mso-position-vertical:absolute;mso-position-vertical-relative:page;mso-width-percent:0;mso-height-percent:0"><v:imagedata r:id="rId7" o:title=" 5300_IBMpos"/><w10:wrap type="square" anchorx="page" anchory="page"/></v:shape></w:pict></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:pict w14:anchorId="5CA4CB64"><v:shape id="Picture 48" o:spid="_x0000_s1026" type="#_x0000_t75" alt="Description: Description: colorblock_PU01" style="position:absolute;margin-left:36pt;margin-top:633.7pt;width:3in;height:89.3pt;z-index:1;visibility:visible;mso-wrap-style:square;mso-wrap-edited:f;mso-width-percent:0;mso-height-percent:0;mso-wrap-distance-left:9.35pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9.35pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:page;mso-position-vertical:absolute;mso-position-vertical-relative:page;mso-width-percent:0;mso-height-percent:0;mso-width-relative:page;mso-height-relative:page"><v:imagedata r:id="rId8" o:title=" colorblock_PU01"/><w10:wrap type="square" anchorx="page" anchory="page"/></v:shape></w:pict></w:r></w:p>
'@

$picRng.InsertXML($picXml)

# ------------------------------------------------------------------
# 2. Remove the "stay informed" social-media blurb at the end of the
#    abstract (IBM Training News / YouTube / Facebook / Twitter),
#    collapsing it down to a single empty paragraph that keeps the
#    intro paragraph's formatting plus the left indent used by the
#    (now removed) bullet lines.
# ------------------------------------------------------------------
$introRng = $d.Content
$foundIntro = $introRng.Find.Execute("To stay informed about IBM training, see the following sites:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundIntro) { throw "Intro paragraph not found" }

$introParaIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $introRng.Start -and $introRng.Start -lt $p.Range.End) {
        $introParaIdx = $i
        break
    }
}
if ($introParaIdx -eq -1) { throw "Could not map intro range to a paragraph" }

$lastRng = $d.Content
$foundLast = $lastRng.Find.Execute("Twitter: twitter.com/websphere_edu", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundLast) { throw "Twitter paragraph not found" }

$lastParaIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $lastRng.Start -and $lastRng.Start -lt $p.Range.End) {
        $lastParaIdx = $i
        break
    }
}
if ($lastParaIdx -eq -1) { throw "Could not map last range to a paragraph" }

# Delete the four "site" paragraphs entirely (including their marks).
$delStart = $d.Paragraphs.Item($introParaIdx + 1).Range.Start
$delEnd = $d.Paragraphs.Item($lastParaIdx).Range.End
$delRng = $d.Range($delStart, $delEnd)
$delRng.Delete()

# Clear the intro paragraph's own run text, keeping its paragraph mark.
$introPara = $d.Paragraphs.Item($introParaIdx)
$textRng = $d.Range($introPara.Range.Start, $introPara.Range.End - 1)
$textRng.Delete()

# It inherits the indentation that used to belong to the bullet lines.
$introPara2 = $d.Paragraphs.Item($introParaIdx)
$introPara2.Format.LeftIndent = 36
